$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 3119447.96
$ws.Range("C9").Value = 324756.99
$ws.Range("D9").Value = 3444204.95
$ws.Range("E9").Value = 9.429084352253774
$ws.Range("F9").Value = 90.57091564774622
$ws.Range("G9").Value = -69.15406245835948
$ws.Range("H9").Value = -63.34724889395527
$ws.Range("I9").Value = -63.98650473110224
$ws.Range("J9").Value = 23100
$ws.Range("K9").Value = 918
$ws.Range("L9").Value = 24018
